# "Avance con la presentacion final"
# Rename the two "Hs Planificadas ..." headers (used by both the sheet's
# table and the line chart's series) to their shorter final names:
#   "Hs Planificadas Totales"     -> "Hs Planificadas"
#   "Hs Planificadas Completadas" -> "Hs Completadas"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cobertura de la Prueba")

# Table header cells (B2:C2). The Excel table ("Tabla1") column names and
# the chart series names are both driven off these two cells, so editing
# the cell text is all that's needed to rename them everywhere they are
# sourced from.
$ws.Range("B2").Value = "Hs Planificadas"
$ws.Range("C2").Value = "Hs Completadas"

# Leave the cursor on D2, matching where the author's selection ended up.
$ws.Range("D2").Select() | Out-Null
